$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 4 and 5 (PHYSICS / MEDICINE study-profile rows) need to swap places
# entirely (all columns A:E), per the commit's reshuffled data.

$a4 = $ws.Range("A4").Value()
$b4 = $ws.Range("B4").Value()
$c4 = $ws.Range("C4").Value()
$d4 = $ws.Range("D4").Value()
$e4 = $ws.Range("E4").Value()

$a5 = $ws.Range("A5").Value()
$b5 = $ws.Range("B5").Value()
$c5 = $ws.Range("C5").Value()
$d5 = $ws.Range("D5").Value()
$e5 = $ws.Range("E5").Value()

$ws.Range("A4").Value = $a5
$ws.Range("B4").Value = $b5
$ws.Range("C4").Value = $c5
$ws.Range("D4").Value = $d5
$ws.Range("E4").Value = $e5

$ws.Range("A5").Value = $a4
$ws.Range("B5").Value = $b4
$ws.Range("C5").Value = $c4
$ws.Range("D5").Value = $d4
$ws.Range("E5").Value = $e4
